# Append a new price-tracking row for 2026-02-07 (Date, Price, Discount, Incredible)
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

$lastRow = $ws.UsedRange.Rows.Count
$newRow = $lastRow + 1

$rng = $ws.Range("A" + $newRow + ":D" + $newRow)

# Force text storage so values such as "2026-02-07" aren't reinterpreted
# as dates/numbers, matching the existing column formatting.
$rng.NumberFormat = "@"

$ws.Range("A" + $newRow).Value = "2026-02-07"
$ws.Range("B" + $newRow).Value = "390000"
$ws.Range("C" + $newRow).Value = "0"
$ws.Range("D" + $newRow).Value = "0"

# Drop the temporary number-format override so the new cells keep the
# same (default/general) style as every other cell in the sheet.
$rng.ClearFormats()
